# Re-order observation records in the "Artfynd" sheet.
#
# Two groups of rows trade places (their Id/TaxonId/species/coords/etc. are
# swapped, but the row number itself stays the position in the sheet):
#   - rows 5 and 6 swap with each other
#   - rows 17, 18 and 19 rotate: new17 = old18, new18 = old19, new19 = old17
#
# Each changed cell is written explicitly (rather than done via a row
# copy/paste) so that cells which only exist in one of the source rows
# (e.g. "M" = Aktivitet, "AC" = Publik kommentar) end up present/absent on
# exactly the right destination row, and cells that don't change (C, J, K,
# L, N, O, S, T, U, V, W, X, Y, AA, AD, AE, AF, AG, AH..., AW, AX, AY, ...)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 <- old row 6
$ws.Range("A5").Value = 130963816
$ws.Range("B5").Value = 79245
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 445932
$ws.Range("R5").Value = 6760079
$ws.Range("Z5").Value = "14:08"
$ws.Range("AB5").Value = "14:08"
$ws.Range("AC5").Value = "Rikligt i närområdet"
$ws.Range("M5").ClearContents()

# Row 6 <- old row 5
$ws.Range("A6").Value = 130960607
$ws.Range("B6").Value = 57884
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 446240
$ws.Range("R6").Value = 6759818
$ws.Range("Z6").Value = "10:26"
$ws.Range("AB6").Value = "10:26"
$ws.Range("AC6").ClearContents()

# Row 17 <- old row 18
$ws.Range("A17").Value = 130961956
$ws.Range("B17").Value = 79864
$ws.Range("E17").Value = 6453
$ws.Range("F17").Value = "Vedskivlav"
$ws.Range("G17").Value = "Hertelidea botryosa"
$ws.Range("H17").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q17").Value = 446084
$ws.Range("R17").Value = 6759981
$ws.Range("AC17").Value = "Miljöbilder"

# Row 18 <- old row 19
$ws.Range("A18").Value = 130960843
$ws.Range("B18").Value = 79245
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("Q18").Value = 446247
$ws.Range("R18").Value = 6759903
$ws.Range("AC18").ClearContents()

# Row 19 <- old row 17
$ws.Range("A19").Value = 130960789
$ws.Range("Q19").Value = 446284
$ws.Range("R19").Value = 6759886
